$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 511.60785
$ws.Range("I15").Value = 511.60785
$ws.Range("K15").Value = 1534.82355
$ws.Range("M15").Value = -1365.82355

$ws.Range("H34").Value = 4220.143
$ws.Range("I34").Value = 4220.143
$ws.Range("K34").Value = 4220.143
$ws.Range("M34").Value = -4017.143

$ws.Range("H36").Value = 4220.143
$ws.Range("I36").Value = 4220.143
$ws.Range("K36").Value = 4220.143
$ws.Range("M36").Value = -3505.143

$ws.Range("H92").Value = 448.33334
$ws.Range("I92").Value = 454.6154
$ws.Range("K92").Value = 454.6154
$ws.Range("M92").Value = 793.3846

$ws.Range("H129").Value = 529.25
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H138").Value = 7657.816
$ws.Range("J138").Value = 7797.1924
$ws.Range("L138").Value = 23391.5772
$ws.Range("N138").Value = -33671.5772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 5972.08
$ws.Range("I6").Value = 2525.2
$ws.Range("J6").Value = 8270
$ws.Range("K6").Value = 2525.2
$ws.Range("L6").Value = 8270
$ws.Range("M6").Value = -2352.2
$ws.Range("N6").Value = -8616

$ws.Range("H24").Value = 16655.666
$ws.Range("J24").Value = 16655.666
$ws.Range("L24").Value = 16655.666
$ws.Range("N24").Value = -17403.666

$ws.Range("H32").Value = 20314.473
$ws.Range("J32").Value = 28939.117
$ws.Range("L32").Value = 28939.117
$ws.Range("N32").Value = -29513.117

$ws.Range("H44").Value = 59999.5
$ws.Range("J44").Value = 59999.5
$ws.Range("L44").Value = 59999.5
$ws.Range("N44").Value = -60975.5

$ws.Range("H45").Value = 2645.8125
$ws.Range("I45").Value = 1791.75
$ws.Range("J45").Value = 3499.875
$ws.Range("K45").Value = 1791.75
$ws.Range("L45").Value = 3499.875
$ws.Range("M45").Value = -1414.75
$ws.Range("N45").Value = -4253.875

$ws.Range("H100").Value = 16655.666
$ws.Range("J100").Value = 16655.666
$ws.Range("L100").Value = 16655.666
$ws.Range("N100").Value = -18819.666

$ws.Range("H122").Value = 718571.94
$ws.Range("I122").Value = 1432285.4
$ws.Range("J122").Value = 4858.4287
$ws.Range("K122").Value = 4296856.199999999
$ws.Range("L122").Value = 14575.2861
$ws.Range("M122").Value = -4294406.199999999
$ws.Range("N122").Value = -19475.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 6544.727
$ws.Range("I107").Value = 4496
$ws.Range("J107").Value = 12008
$ws.Range("K107").Value = 4496
$ws.Range("L107").Value = 12008
$ws.Range("M107").Value = -2576
$ws.Range("N107").Value = -15848

$ws.Range("H134").Value = 3075.2173
$ws.Range("I134").Value = 1534.4615
$ws.Range("K134").Value = 4603.3845
$ws.Range("M134").Value = -2068.3845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H31").Value = 4549.4287
$ws.Range("I31").Value = 3240.7334
$ws.Range("K31").Value = 3240.7334
$ws.Range("M31").Value = -2945.7334

$ws.Range("H34").Value = 4549.4287
$ws.Range("I34").Value = 3240.7334
$ws.Range("K34").Value = 3240.7334
$ws.Range("M34").Value = -3038.7334

$ws.Range("H42").Value = 9500
$ws.Range("J42").Value = 9500
$ws.Range("L42").Value = 9500
$ws.Range("N42").Value = -10686

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H94").Value = 1006
$ws.Range("J94").Value = 1099.8
$ws.Range("L94").Value = 1099.8
$ws.Range("N94").Value = -2001.8

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H107").Value = 420.55554
$ws.Range("J107").Value = 454.83334
$ws.Range("L107").Value = 454.83334
$ws.Range("N107").Value = -4294.83334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 1225
$ws.Range("J52").Value = 1225
$ws.Range("L52").Value = 3675
$ws.Range("N52").Value = -4207

$ws.Range("H80").Value = 2550
$ws.Range("I80").Value = 3000.3333
$ws.Range("K80").Value = 9000.999899999999
$ws.Range("M80").Value = -8064.999899999999

$ws.Range("H83").Value = 2550
$ws.Range("I83").Value = 3000.3333
$ws.Range("K83").Value = 27002.9997
$ws.Range("M83").Value = -22322.9997

$ws.Range("H139").Value = 9259.5
$ws.Range("J139").Value = 11249
$ws.Range("L139").Value = 33747
$ws.Range("N139").Value = -44027

$ws.Range("H140").Value = 5023
$ws.Range("I140").Value = 4364
$ws.Range("K140").Value = 13092
$ws.Range("M140").Value = -7912

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 34999
$ws.Range("J39").Value = 34999
$ws.Range("L39").Value = 34999
$ws.Range("N39").Value = -36063

$ws.Range("H80").Value = 9281.6
$ws.Range("J80").Value = 9240.25
$ws.Range("L80").Value = 9240.25
$ws.Range("N80").Value = -11236.25

$ws.Range("H83").Value = 9281.6
$ws.Range("J83").Value = 9240.25
$ws.Range("K83").Value = 46201.25
$ws.Range("N83").Value = -56185.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2873747.8
$ws.Range("J2").Value = 21246.6
$ws.Range("L2").Value = 21246.6
$ws.Range("N2").Value = -21470.6

$ws.Range("H14").Value = 6761
$ws.Range("J14").Value = 6761
$ws.Range("L14").Value = 6761
$ws.Range("N14").Value = -7105

$ws.Range("H46").Value = 2909
$ws.Range("I46").Value = 2250
$ws.Range("J46").Value = 3285.5715
$ws.Range("K46").Value = 2250
$ws.Range("L46").Value = 3285.5715
$ws.Range("M46").Value = -2062
$ws.Range("N46").Value = -3661.5715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 11968.667
$ws.Range("J5").Value = 11968.667
$ws.Range("L5").Value = 11968.667
$ws.Range("N5").Value = -12192.667

$ws.Range("H54").Value = 285
$ws.Range("I54").Value = 285
$ws.Range("K54").Value = 285
$ws.Range("M54").Value = 235

$ws.Range("H100").Value = 2595.8572
$ws.Range("I100").Value = 3338.111
$ws.Range("J100").Value = 1259.8
$ws.Range("K100").Value = 6676.222
$ws.Range("L100").Value = 2519.6
$ws.Range("M100").Value = -6135.222
$ws.Range("N100").Value = -3601.6

$ws.Range("H107").Value = 1882.7142
$ws.Range("I107").Value = 743
$ws.Range("K107").Value = 2229
$ws.Range("M107").Value = -309

$ws.Range("H113").Value = 1122.3846
$ws.Range("I113").Value = 673.5455
$ws.Range("J113").Value = 1451.5333
$ws.Range("K113").Value = 2020.6365
$ws.Range("L113").Value = 4354.5999
$ws.Range("M113").Value = 149.3635000000002
$ws.Range("N113").Value = -8694.599900000001

$ws.Range("H132").Value = 1556
$ws.Range("I132").Value = 663.46155
$ws.Range("K132").Value = 1990.38465
$ws.Range("M132").Value = 539.61535
